# The "WALL CABINETS" section header (previously mislabeled "2 DOORS")
# needs to be corrected for rows 27-41 in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$ws.Range("A27:A41").Value = "WALL CABINETS"

# Reflect the author's on-screen selection/scroll state at save time:
# the range A27:A41 is selected with A27 as the active cell, and the
# view is scrolled down so row 141 is visible at the top.
$excel.ActiveWindow.ScrollRow = 141
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A27:A41").Select()
